$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting rows 24:30 down to 25:31
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new weekly data point
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C24").Value = 'Ñuble'
$ws.Range("D24").Value = 44508
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112013
$ws.Range("G24").Value = 'Alcachofa'
$ws.Range("H24").Value = 'Madrigal'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 11000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 11500
$ws.Range("N24").Value = '$/caja 40 unidades'
$ws.Range("O24").Value = 'Provincia del Elquí'
$ws.Range("P24").Value = 288
$ws.Range("Q24").Value = 40
$ws.Range("R24").Value = 'Hortaliza'
